$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D27").Value = "개발자를 위한 AWS 클라우드 보안 (2) - 로깅 및 모니터링과 데이터 보호"
$ws.Range("E27").Value = "https://blog.pingpong.us/aws-cloud-security-for-devs-2/"

$ws.Range("D32").Value = "ICE (Individual conditional expectation)"

$ws.Range("D42").Value = "파이썬 모듈 설치파일 로컬 다운로드"
$ws.Range("E42").Value = "https://kjk92.tistory.com/84"

$ws.Range("D46").Value = "[한국생명공학연구원] 2022년 06월, 생물정보학(Bioinformatics 채용), 유전체 연구직"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/476"
